# Fix 2016 social services immigration office budget figures.
# The budget figures in column E (rows 2-40) were mistakenly entered in
# plain Shekels instead of thousands of Shekels (as labelled by the
# "תקציב 2016 אלפי ₪" header), so every raw input value is divided by
# 1,000. The per-section subtotal cells (rows 5, 8, 14, 21, 24, 27, 38)
# held SUM() formulas over those inputs; they are replaced with their
# corrected plain values (matching what a Paste-Special "Divide" over
# E2:E40 would do to formula cells). The grand-subtotal cells E41 and
# E42 keep their original SUM formulas and simply recalculate using the
# corrected inputs. The now-unused helper numbers that used to live in
# column J (rows 26-34, one of which held a SUM formula) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section: אגף תעסוקה (rows 2-4, subtotal row 5) ---
$ws.Range("E2").Value = 62634
$ws.Range("E3").Value = 6582
$ws.Range("E4").Value = 4500
$ws.Range("E5").Value = 73716

# --- Section: קליטה במדע (rows 6-7, subtotal row 8) ---
$ws.Range("E6").Value = 12000
$ws.Range("E7").Value = 5800
$ws.Range("E8").Value = 17800

# --- Section: אגף בכיר לקליטה בקהילה (rows 9-13, subtotal row 14) ---
$ws.Range("E9").Value = 14000
$ws.Range("E10").Value = 2150
$ws.Range("E11").Value = 420
$ws.Range("E12").Value = 300
$ws.Range("E13").Value = 381
$ws.Range("E14").Value = 17251

# --- Section: אגף צעירים (rows 15-20, subtotal row 21) ---
$ws.Range("E15").Value = 5913.2
$ws.Range("E16").Value = 2754.12
$ws.Range("E17").Value = 619
$ws.Range("E18").Value = 644
$ws.Range("E19").Value = 1600
$ws.Range("E20").Value = 300
$ws.Range("E21").Value = 11830.32

# --- Section: יזמות (rows 22-23, subtotal row 24) ---
$ws.Range("E22").Value = 6800
$ws.Range("E23").Value = 1800
$ws.Range("E24").Value = 8600

# --- Section: עידוד עליה (rows 25-26, subtotal row 27) ---
$ws.Range("E25").Value = 4500
$ws.Range("E26").Value = 3500
$ws.Range("E27").Value = 8000

# --- Section: שירות הרווחה (rows 28-37, subtotal row 38) ---
$ws.Range("E28").Value = 2528
$ws.Range("E29").Value = 1374.24
$ws.Range("E30").Value = 594
$ws.Range("E31").Value = 562.5
$ws.Range("E32").Value = 314.1
$ws.Range("E33").Value = 503.1
$ws.Range("E34").Value = 503
$ws.Range("E35").Value = 231.525
$ws.Range("E36").Value = 200
$ws.Range("E37").Value = 182.745
$ws.Range("E38").Value = 6993.21

# --- Section: מנהל דיור (rows 39-40, subtotal row 41 keeps its formula) ---
$ws.Range("E39").Value = 1250
$ws.Range("E40").Value = 1000

# Clear the stray helper numbers (and their SUM formula) that used to
# live in column J - they duplicated the (wrong) subtotal figures.
$ws.Range("J26:J34").ClearContents()

# Recalculate so the remaining formulas (E41, E42, and the grand total)
# pick up the corrected figures.
$excel.Calculate()

# Restore the frozen-pane scroll position (the view had drifted down to
# row 32) and leave the selection on the now-empty helper column.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
[void]$ws.Range("G1:J1048576").Select()
